$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.630.21'
$ws.Range("E2").Value = '  -1.55%  '

# Row 3
$ws.Range("D3").Value = '3.413.26'
$ws.Range("E3").Value = '  -1.94%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").Value = '''597.01'
$ws.Range("E5").Value = '  -0.89%  '

# Row 6
$ws.Range("D6").Value = '''142.08'
$ws.Range("E6").Value = '  -4.00%  '

# Row 7
$ws.Range("D7").Value = '3.412.54'
$ws.Range("E7").Value = '  -1.86%  '

# Row 8
$ws.Range("E8").Value = '  -0.13%  '

# Row 9
$ws.Range("D9").Value = '''0.470'
$ws.Range("E9").Value = '  -2.48%  '

# Row 10
$ws.Range("D10").Value = '''7.94'
$ws.Range("E10").Value = '  +4.94%  '

# Row 11
$ws.Range("D11").Value = '''0.134'
$ws.Range("E11").Value = '  -5.80%  '

# Row 12
$ws.Range("D12").Value = '''0.405'
$ws.Range("E12").Value = '  -4.33%  '

# Row 13
$ws.Range("D13").Value = '3.985.88'
$ws.Range("E13").Value = '  -2.07%  '

# Row 14
$ws.Range("D14").Value = '''0.0000200'
$ws.Range("E14").Value = '  -6.14%  '

# Row 15
$ws.Range("D15").Value = '''29.55'
$ws.Range("E15").Value = '  -6.05%  '

# Row 16
$ws.Range("E16").Value = '  -0.67%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.405.69'
$ws.Range("E17").Value = '  -2.91%  '

# Row 18
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '65.616.63'
$ws.Range("E18").Value = '  -1.68%  '

# Row 19
$ws.Range("D19").Value = '''10.35'
$ws.Range("E19").Value = '  +2.90%  '

# Row 20
$ws.Range("D20").Value = '''6.11'
$ws.Range("E20").Value = '  -5.23%  '

# Row 21
$ws.Range("D21").Value = '''14.52'
$ws.Range("E21").Value = '  -5.41%  '

# Row 22
$ws.Range("D22").Value = '''414.39'
$ws.Range("E22").Value = '  -5.42%  '

# Row 23
$ws.Range("D23").Value = '''0.577'
$ws.Range("E23").Value = '  -5.19%  '

# Row 24
$ws.Range("D24").Value = '''77.15'
$ws.Range("E24").Value = '  -2.90%  '

# Row 25
$ws.Range("E25").Value = '  +0.09%  '

# Row 26
$ws.Range("D26").Value = '3.545.21'
$ws.Range("E26").Value = '  -2.07%  '

# Row 27
$ws.Range("D27").Value = '''0.0000109'
$ws.Range("E27").Value = '  -8.65%  '

# Row 28
$ws.Range("D28").Value = '''9.24'
$ws.Range("E28").Value = '  -5.79%  '

# Row 29
$ws.Range("D29").Value = '''7.80'
$ws.Range("E29").Value = '  -6.73%  '

# Row 30
$ws.Range("D30").Value = '''2.42'
$ws.Range("E30").Value = '  -2.70%  '

# Row 31
$ws.Range("D31").Value = '''0.999'
$ws.Range("E31").Value = '  -0.37%  '

# Row 32
$ws.Range("E32").Value = '  -5.01%  '

# Row 33
$ws.Range("D33").Value = '''1.45'
$ws.Range("E33").Value = '  -8.24%  '

# Row 34
$ws.Range("D34").Value = '''24.50'
$ws.Range("E34").Value = '  -3.46%  '

# Row 35
$ws.Range("D35").Value = '3.405.49'
$ws.Range("E35").Value = '  -1.93%  '

# Row 37
$ws.Range("D37").Value = '''1.68'
$ws.Range("E37").Value = '  -6.68%  '

# Row 38
$ws.Range("E38").Value = '  -8.91%  '

# Row 39
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '''7.50'
$ws.Range("E39").Value = '  -5.37%  '

# Row 40
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '''0.999'
$ws.Range("E40").Value = '  -0.01%  '

# Row 41
$ws.Range("D41").Value = '''168.90'
$ws.Range("E41").Value = '  -4.47%  '

# Row 42
$ws.Range("D42").Value = '''0.0852'
$ws.Range("E42").Value = '  -3.72%  '

# Row 43
$ws.Range("D43").Value = '''0.871'
$ws.Range("E43").Value = '  -2.34%  '

# Row 44
$ws.Range("D44").Value = '''5.02'
$ws.Range("E44").Value = '  -7.40%  '

# Row 45
$ws.Range("D45").Value = '''1.89'
$ws.Range("E45").Value = '  -10.98%  '

# Row 46
$ws.Range("D46").Value = '''45.36'
$ws.Range("E46").Value = '  -2.02%  '

# Row 47
$ws.Range("D47").Value = '''26.23'
$ws.Range("E47").Value = '  -8.78%  '

# Row 48
$ws.Range("E48").Value = '  -3.93%  '

# Row 49
$ws.Range("D49").Value = '''7.03'
$ws.Range("E49").Value = '  -5.60%  '

# Row 50
$ws.Range("E50").Value = '  -6.89%  '

# Row 51
$ws.Range("D51").Value = '''0.914'
$ws.Range("E51").Value = '  -6.68%  '
